# Weitere Vorbereitungen für den 2015er launch
# Remove the "Nicki Westhoff" player row (row 6) from the player sheet;
# subsequent rows shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(6).Delete()
